# Generate Report for Handback
# Update the "latest generated" timestamps for the second tracked file
# (7fa5d33f-0b67-4675-93f3-dc92d1eec877.md) after a new handback round.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# zh-cn: Correspond Handoff Datetime / Correspond Handback DateTime for row 3
$wsZhCn.Range("H3").Value = "2016-08-25 10:50:27"
$wsZhCn.Range("K3").Value = "2016-08-25 10:50:43"

# de-de: Correspond Handoff Datetime / Correspond Handback DateTime for row 3
$wsDeDe.Range("H3").Value = "2016-08-25 10:50:32"
$wsDeDe.Range("K3").Value = "2016-08-25 10:50:51"

# Overview: Latest HO Xliff Generate Date mirrors the de-de handoff datetime
$wsOverview.Range("G3").Value = "2016-08-25 10:50:32"
